$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "Scalar_annualized1"
$ws.Range("B37").Value = "Test scalar annulized for value"
$ws.Range("C37").Value = "scalar_annualized_test1"

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("C40").Select()
